$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 109; this shifts the existing rows 109-223 down to 110-224.
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with the new weekly price record.
$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44781
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = 100112017
$ws.Range("G109").Value = "Apio"
$ws.Range("H109").Value = "Americana (o)"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 120
$ws.Range("K109").Value = 8500
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = 8750
$ws.Range("N109").Value = "$/docena de matas"
$ws.Range("O109").Value = "Provincia del Elquí"
$ws.Range("P109").Value = 1458
$ws.Range("Q109").Value = 6
$ws.Range("R109").Value = "Hortaliza"
